# ---------------------------------------------------------------------------
# landings_fabric_vessels_per_species_region_month_2017.xlsx
# "running script including 2017 data"
#
# The meaningful, content-level edit in this commit is a worksheet rename:
#   "bf_nac_región" -> "bf_nac_region"
# (dropping the accented/non-ASCII character so the tab name is safe to
# reference from the ingestion script), together with the defined names
# ( _xlnm._FilterDatabase / _xlnm.Print_Titles ) that point at it.
#
# The formula blocks in the totals rows (23, 24, 27) of "bf_nac_region" are
# re-entered as single fill-across ranges, which is how this workbook's
# totals rows were originally produced by the reporting script - Excel
# collapses same-pattern formulas dragged across a row into a shared-formula
# group when the sheet is saved.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the national/region sheet -----------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "bf_nac_region"

# --- 2. Repoint the defined names so they keep referring to the renamed
#        sheet (Excel auto-updates most references on rename, but make sure
#        every name's RefersTo is explicit/correct) -------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "bf_nac_region!_FilterDatabase") {
        $n.RefersTo = "=bf_nac_region!`$A`$5:`$P`$18"
    }
    if ($n.Name -eq "bf_nac_region!Print_Titles") {
        $n.RefersTo = "=bf_nac_region!`$1:`$5"
    }
    if ($n.Name -eq "bf_nac_mes!_FilterDatabase") {
        $n.RefersTo = "=bf_nac_mes!#REF!"
    }
    if ($n.Name -eq "bf_nac_mes!Print_Titles") {
        $n.RefersTo = "=bf_nac_mes!`$1:`$5"
    }
}

# --- 3. Re-enter the totals-row formulas as row fills so they serialize as
#        shared-formula groups, same as the source workbook ----------------
$ws1.Range("B23:P23").Formula = "=SUM(B6:B18)"
$ws1.Range("B24").Formula = "=SUM(B20)"
$ws1.Range("C24:P24").Formula = "=SUM(C20)"
$ws1.Range("B27").Formula = "=SUM(B22:B26)"
$ws1.Range("C27:P27").Formula = "=SUM(C22:C26)"
